# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E22) listed the mora periods in descending
# order (1701, 1612, 1611, 1610, 1609, 1608, 1607). This update re-sorts
# them in ascending order (1607, 1608, 1609, 1610, 1611, 1612, 1701) to
# reflect the refreshed EC database / newly appended statement period.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periodos = @("1607", "1608", "1609", "1610", "1611", "1612", "1701")

for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periodos[$i]
}
